$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.249.53'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '1.906.40'
$ws.Range("E3").Value = '  +0.32%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.51%  '

$ws.Range("D5").Value = '''307.44'
$ws.Range("E5").Value = '  +1.21%  '

$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").Value = '''0.5263'
$ws.Range("E7").Value = '  +1.33%  '

$ws.Range("D8").Value = '''0.3816'
$ws.Range("E8").Value = '  +1.45%  '

$ws.Range("D9").Value = '''0.07293'
$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("D10").Value = '''22.01'
$ws.Range("E10").Value = '  +4.10%  '

$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("D12").Value = '''0.08173'
$ws.Range("E12").Value = '  -2.35%  '

$ws.Range("D13").Value = '''95.89'
$ws.Range("E13").Value = '  -0.97%  '

$ws.Range("D14").Value = '''5.358'
$ws.Range("E14").Value = '  +1.37%  '

$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = '''1.001'
$ws.Range("E15").Value = '  +0.47%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.000008644'
$ws.Range("E16").Value = '  +0.08%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '''14.76'
$ws.Range("E17").Value = '  +1.69%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '1.393.64'
$ws.Range("E18").Value = '  -26.69%  '

$ws.Range("D19").Value = '''1.002'
$ws.Range("E19").Value = '  +0.28%  '

$ws.Range("D20").Value = '27.289.29'
$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").Value = '''5.093'
$ws.Range("E21").Value = '  +0.31%  '

$ws.Range("D22").Value = '''10.81'
$ws.Range("E22").Value = '  +1.66%  '

$ws.Range("D23").Value = '''6.525'
$ws.Range("E23").Value = '  +1.40%  '

$ws.Range("D24").Value = '''149.66'
$ws.Range("E24").Value = '  +2.31%  '

$ws.Range("D25").Value = '''2.305'
$ws.Range("E25").Value = '  -1.17%  '

$ws.Range("E26").Value = '  +0.35%  '

$ws.Range("D27").Value = '''1.735'
$ws.Range("E27").Value = '  -0.56%  '

$ws.Range("D28").Value = '''116.67'
$ws.Range("E28").Value = '  +1.51%  '

$ws.Range("D29").Value = '''4.849'
$ws.Range("E29").Value = '  +0.71%  '

$ws.Range("D30").Value = '''4.825'
$ws.Range("E30").Value = '  -1.24%  '

$ws.Range("D31").Value = '''0.09257'
$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("D32").Value = '''0.8291'
$ws.Range("E32").Value = '  +4.00%  '

$ws.Range("D33").Value = '''0.05076'
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("D34").Value = '''1.228'
$ws.Range("E34").Value = '  -1.06%  '

$ws.Range("D35").Value = '''2.998'
$ws.Range("E35").Value = '  +1.44%  '

$ws.Range("D36").Value = '''3.355'
$ws.Range("E36").Value = '  -1.65%  '

$ws.Range("D37").Value = '''2.686'
$ws.Range("E37").Value = '  +3.30%  '

$ws.Range("D38").Value = '''0.5801'
$ws.Range("E38").Value = '  +2.27%  '

$ws.Range("D39").Value = '''0.01999'
$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("D40").Value = '''1.080'
$ws.Range("E40").Value = '  +0.54%  '

$ws.Range("D41").Value = '''9.262'
$ws.Range("E41").Value = '  +2.52%  '

$ws.Range("D42").Value = '''6.552'
$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").Value = '''116.54'
$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("D44").Value = '''0.1522'
$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").Value = '''0.4922'
$ws.Range("E45").Value = '  +1.53%  '

$ws.Range("E46").Value = '  +0.61%  '

$ws.Range("D47").Value = '''1.001'
$ws.Range("E47").Value = '  +0.43%  '

$ws.Range("D48").Value = '''1.639'
$ws.Range("E48").Value = '  +0.38%  '

$ws.Range("D49").Value = '''38.76'
$ws.Range("E49").Value = '  +2.71%  '

$ws.Range("D50").Value = '''0.06199'
$ws.Range("E50").Value = '  +4.25%  '

$ws.Range("D51").Value = '''64.18'
$ws.Range("E51").Value = '  +0.42%  '
